# edit.ps1 - applies the "Added email to writeup" commit to Writeup.docx
#
# Strategy: use Range.InsertXML() with full <w:p>...</w:p> replacement
# blocks (whole-paragraph granularity) for every paragraph whose content
# changes. InsertXML is unreliable when fed "loose" runs that don't cover
# an entire paragraph, so every call below targets exactly one (or, for
# the DinnerTable/stdout split, exactly two) complete paragraphs at once.
#
# Order of operations: first perform all the edits that do not change the
# total paragraph count (so Paragraphs(n) indices stay valid throughout),
# then finish with the structural insert (new email paragraph) last.

$d = $word.ActiveDocument
$wns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- 1. "Homework 1 Writeup" -> split off "Writeup" as a flagged spelling word ---
$pTitle = $d.Paragraphs(3)
if ($pTitle.Range.Text.TrimEnd([char]13,[char]7) -ne "Homework 1 Writeup") {
    throw "paragraph 3 mismatch: $($pTitle.Range.Text)"
}
$xmlTitle = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">Homework 1 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Writeup</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$pTitle.Range.InsertXML($xmlTitle)

# --- 2. "...a lot of playtesting, the fact..." -> flag "playtesting" ---
$pPlay = $d.Paragraphs(9)
$xmlPlay = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">I ran into various pitfalls, or rather plateaus, but it turns out that after a lot of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>playtesting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, the fact of sitting </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space="preserve"> unhappiest people at the corners increase</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> my results by a good 20%. The three files t</w:t></w:r><w:r><w:t>ake under</w:t></w:r><w:r><w:t xml:space="preserve"> 30 seconds of runtime.</w:t></w:r></w:p>
'@
$pPlay.Range.InsertXML($xmlPlay)

# --- 3. "java DinnerTable "input.txt"" / "I had routed stdout..." paragraphs ---
#        split across both paragraphs together (DinnerTable + stdout flagged)
$pCode = $d.Paragraphs(13)
$pStdout = $d.Paragraphs(14)
$rngCode = $d.Range($pCode.Range.Start, $pStdout.Range.End)
$xmlCode = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">To run the code with an input file: </w:t></w:r><w:r w:rsidRPr="001271DC"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve">java </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>DinnerTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve"> "input.txt"</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">I had routed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stdout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to file, see: hw1-soln1.txt, hw1-soln2.txt and hw1-soln3.txt</w:t></w:r></w:p>
'@
$rngCode.InsertXML($xmlCode)

# --- 4. "Java with notepad++, and git for source code management." -> wrap whole
#        sentence as a flagged grammar run (gramStart before 1st run, gramEnd after last) ---
$pSoftware = $d.Paragraphs(17)
$xmlSoftware = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Java with notepad++, </w:t></w:r><w:r><w:t>and git for source code</w:t></w:r><w:r><w:t xml:space="preserve"> management.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$pSoftware.Range.InsertXML($xmlSoftware)

# --- 5. "Home computer, 3.1ghz quad core..." -> flag "3.1ghz" as a grammar item ---
$pHardware = $d.Paragraphs(20)
$xmlHardware = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Home computer, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>3.1ghz</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> quad core i5 processor, 8 GB's memory @1600.</w:t></w:r></w:p>
'@
$pHardware.Range.InsertXML($xmlHardware)

# --- 6. Insert the new right-justified email paragraph right after "Andrew Helenius"
#        (and before the "10/20/2013" date paragraph). Done last since it changes
#        the paragraph count and would shift every index used above. ---
$pName = $d.Paragraphs(1)
$pName.Range.InsertParagraphAfter()
$pEmail = $d.Paragraphs(2)
$pEmail.Range.Text = "andrew.helenius@gmail.com"

Write-Output "done"
